# Apply updated cryptocurrency price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value. Cells whose new text would be
# auto-parsed as a number (losing trailing zeros / leading zero formatting like
# "0.0310") are forced to Text format first so the stored value matches exactly.
$updates = @(
    @{ Cell = 'D2'; Value = '42.959.41'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -5.98%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '2.543.24'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -1.92%  '; ForceText = $false }
    @{ Cell = 'D4'; Value = '0.999'; ForceText = $true }
    @{ Cell = 'E4'; Value = '  -0.10%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '299.34'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -2.78%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '94.19'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -4.35%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  -3.16%  '; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.548'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -4.53%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '35.99'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -6.39%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.0805'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -3.80%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '7.73'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -3.54%  '; ForceText = $false }
    @{ Cell = 'E13'; Value = '  +3.57%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '2.937.00'; ForceText = $false }
    @{ Cell = 'E14'; Value = '  -1.90%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '2.560.49'; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -1.52%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '0.870'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -3.96%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '14.08'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -4.16%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '43.006.54'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  -6.04%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '13.08'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +3.85%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.0₃0977'; ForceText = $false }
    @{ Cell = 'E20'; Value = '  -3.01%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '6.60'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -0.79%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '71.75'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -2.11%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '257.59'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -9.29%  '; ForceText = $false }
    @{ Cell = 'E24'; Value = '  -2.84%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '2.13'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -4.74%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '29.04'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +0.27%  '; ForceText = $false }
    @{ Cell = 'E27'; Value = '  +0.09%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  -5.18%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '37.61'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -1.65%  '; ForceText = $false }
    @{ Cell = 'E30'; Value = '  -3.32%  '; ForceText = $false }
    @{ Cell = 'E31'; Value = '  -4.30%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '153.79'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -2.32%  '; ForceText = $false }
    @{ Cell = 'E33'; Value = '  -1.61%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '2.16'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -3.44%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '3.37'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -6.07%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '0.0795'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -3.92%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  -5.17%  '; ForceText = $false }
    @{ Cell = 'E38'; Value = '  -2.24%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '17.08'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +7.99%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '23.41'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +9.50%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '3.43'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -2.13%  '; ForceText = $false }
    @{ Cell = 'B42'; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false }
    @{ Cell = 'D42'; Value = '3.89'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -2.53%  '; ForceText = $false }
    @{ Cell = 'B43'; Value = 'VeChain'; ForceText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.0310'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -4.47%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '2.076.65'; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -1.07%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +0.06%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '84.62'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -9.76%  '; ForceText = $false }
    @{ Cell = 'E47'; Value = '  -3.29%  '; ForceText = $false }
    @{ Cell = 'E48'; Value = '  +1.94%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '2.794.35'; ForceText = $false }
    @{ Cell = 'E49'; Value = '  -1.86%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '104.38'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -3.53%  '; ForceText = $false }
    @{ Cell = 'E51'; Value = '  -3.56%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
